# Update "想去人数" (interest count) values in the F column of the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets, reflecting
# a refreshed data pull (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# --- "展览" sheet (sheet1) ---
$wsExhibition.Range("F2").Value  = 14976
$wsExhibition.Range("F3").Value  = 18904
$wsExhibition.Range("F5").Value  = 133
$wsExhibition.Range("F13").Value = 56
$wsExhibition.Range("F14").Value = 145
$wsExhibition.Range("F15").Value = 215
$wsExhibition.Range("F17").Value = 1445
$wsExhibition.Range("F20").Value = 94
$wsExhibition.Range("F22").Value = 7868
$wsExhibition.Range("F23").Value = 990
$wsExhibition.Range("F24").Value = 31
$wsExhibition.Range("F27").Value = 1236
$wsExhibition.Range("F29").Value = 6025
$wsExhibition.Range("F32").Value = 167
$wsExhibition.Range("F34").Value = 275
$wsExhibition.Range("F35").Value = 5390
$wsExhibition.Range("F37").Value = 5
$wsExhibition.Range("F39").Value = 45

# --- "全部类型" sheet (sheet4) ---
$wsAllTypes.Range("F2").Value  = 14976
$wsAllTypes.Range("F3").Value  = 18904
$wsAllTypes.Range("F5").Value  = 133
$wsAllTypes.Range("F13").Value = 56
$wsAllTypes.Range("F14").Value = 145
$wsAllTypes.Range("F15").Value = 215
$wsAllTypes.Range("F17").Value = 1445
$wsAllTypes.Range("F21").Value = 94
$wsAllTypes.Range("F23").Value = 7868
$wsAllTypes.Range("F24").Value = 990
$wsAllTypes.Range("F25").Value = 31
$wsAllTypes.Range("F28").Value = 1236
$wsAllTypes.Range("F32").Value = 6025
$wsAllTypes.Range("F35").Value = 167
$wsAllTypes.Range("F37").Value = 275
$wsAllTypes.Range("F38").Value = 5390
$wsAllTypes.Range("F40").Value = 5
$wsAllTypes.Range("F42").Value = 45
